$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 83, shifting existing rows 83-193 down to 84-194
$ws.Rows.Item(83).Insert()

$ws.Cells.Item(83, 1).Value = 10
$ws.Cells.Item(83, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(83, 3).Value = "La Araucanía"
$ws.Cells.Item(83, 4).Value = 45033
$ws.Cells.Item(83, 5).Value = 9
$ws.Cells.Item(83, 6).Value = 100114002
$ws.Cells.Item(83, 7).Value = "Camote"
$ws.Cells.Item(83, 8).Value = "Sin especificar"
$ws.Cells.Item(83, 9).Value = "Primera"
$ws.Cells.Item(83, 10).Value = 25
$ws.Cells.Item(83, 11).Value = 26000
$ws.Cells.Item(83, 12).Value = 26000
$ws.Cells.Item(83, 13).Value = 26000
$ws.Cells.Item(83, 14).Value = '$/caja 18 kilos'
$ws.Cells.Item(83, 15).Value = "Perú"
$ws.Cells.Item(83, 16).Value = 1444
$ws.Cells.Item(83, 17).Value = 18
$ws.Cells.Item(83, 18).Value = "Hortaliza"
